$d = $word.ActiveDocument

# Explicitly set PageBreakBefore = False on every paragraph in the
# document body (mirrors the diff adding <w:pageBreakBefore w:val="0"/>
# to each <w:p>/<w:pPr> in document.xml).
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = $false
}

# Also set PageBreakBefore = False on the built-in heading / title /
# subtitle paragraph styles (mirrors the diff adding the same element to
# each of those styles' <w:pPr> in styles.xml).
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", `
                "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $s.ParagraphFormat.PageBreakBefore = $false
}

Write-Output "done"
